# Adds a new "comment" column to every transaction sheet in the workbook.
# (author's commit: "forgot to commit pieces, commit local progress now")
#
# xlPasteFormats = -4122 ; used so the new header cell inherits the same
# cell style as its left-hand neighbour (font/bold/etc.), the same way the
# existing header cells in each sheet already do.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# rsu: F1 ("currency") -> add G1 ("comment")
$ws = $wb.Worksheets.Item("rsu")
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial($xlPasteFormats)
$ws.Range("G1").Value = "comment"

# espp: F1 ("currency") -> add G1 ("comment")
$ws = $wb.Worksheets.Item("espp")
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial($xlPasteFormats)
$ws.Range("G1").Value = "comment"

# dividends: E1 ("currency") -> add F1 ("comment"); G1 gets the same
# (empty) style as a trailing placeholder cell, matching buy_orders.
$ws = $wb.Worksheets.Item("dividends")
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial($xlPasteFormats)
$ws.Range("F1").Value = "comment"
$ws.Range("E1").Copy()
$ws.Range("G1").PasteSpecial($xlPasteFormats)

# buy_orders: G1 already exists as a styled-but-empty placeholder cell;
# it just needs the "comment" label written into it.
$ws = $wb.Worksheets.Item("buy_orders")
$ws.Range("G1").Value = "comment"

# sell_orders: F1 ("currency") -> add G1 ("comment")
$ws = $wb.Worksheets.Item("sell_orders")
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial($xlPasteFormats)
$ws.Range("G1").Value = "comment"

# currency_conversions: E1 ("tax_withholding") -> add F1 ("comment")
$ws = $wb.Worksheets.Item("currency_conversions")
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial($xlPasteFormats)
$ws.Range("F1").Value = "comment"
